$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.099.54"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.316.78"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.28"
$ws.Range("E5").Value = "  -5.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.25"
$ws.Range("E6").Value = "  +4.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.36"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0915"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.29"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.980"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.667.73"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.317.33"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.202.60"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  -4.44%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.73"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("E22").Value = "  -8.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.13"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.30"
$ws.Range("E25").Value = "  -8.72%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.83"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.72"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.16"
$ws.Range("E32").Value = "  -7.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").Value = "  -5.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.87"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("E35").Value = "  +12.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.130"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("E39").Value = "  -6.29%  "
$ws.Range("E40").Value = "  -4.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.53"
$ws.Range("E41").Value = "  +7.68%  "
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.87"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.25"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.90"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.01"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.88"
$ws.Range("E50").Value = "  +5.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("E51").Value = "  -0.76%  "
